$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
# Row 16 and Row 18 swap their "Periodo Mora" (E) / "Valor Mora" (F) values;
# Row 17 stays as-is.
$ws.Range("E16").Value = "1809"
$ws.Range("F16").Value = 31249

$ws.Range("E18").Value = "1902"
$ws.Range("F18").Value = 28124
